# Add Infer to correlation
# 1) Update existing correlation numbers on all_tools / checker_framework / typestate_checker
# 2) Add a new "infer" worksheet (after typestate_checker) with its own correlation table

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the re-computed correlation stats on the three existing sheets
# ---------------------------------------------------------------------------

$wsAllTools = $wb.Worksheets.Item("all_tools")
$wsAllTools.Range("D3").Value = 378
$wsAllTools.Range("F3").Value = -0.03275384325686242
$wsAllTools.Range("G3").Value = 0.6458960401104946
$wsAllTools.Range("H3").Value = -0.04216656825659747
$wsAllTools.Range("I3").Value = 0.6770078435983685

$wsChecker = $wb.Worksheets.Item("checker_framework")
$wsChecker.Range("F3").Value = -0.1855044244235976
$wsChecker.Range("G3").Value = 0.02021255303764821
$wsChecker.Range("H3").Value = -0.2319964250630901
$wsChecker.Range("I3").Value = 0.02020171421364228

$wsTypestate = $wb.Worksheets.Item("typestate_checker")
$wsTypestate.Range("D3").Value = 311
$wsTypestate.Range("F3").Value = 0.002808180844581953
$wsTypestate.Range("G3").Value = 0.9687827959561681
$wsTypestate.Range("H3").Value = 0.006327541624453262
$wsTypestate.Range("I3").Value = 0.9501801871045381

# ---------------------------------------------------------------------------
# 2. Add the new "infer" worksheet after typestate_checker
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsInfer = $wb.Worksheets.Add($null, $lastSheet)
$wsInfer.Name = "infer"

# Header row
$wsInfer.Range("A1").Value = "Complexity Metric"
$wsInfer.Range("B1").Value = "# of snippets judged (complexity)"
$wsInfer.Range("C1").Value = "# of snippets with warnings"
$wsInfer.Range("D1").Value = "# of warnings"
$wsInfer.Range("E1").Value = "# of data points for correlation"
$wsInfer.Range("F1").Value = "Kendall's Tau (" + [char]0x03C4 + ")"
$wsInfer.Range("G1").Value = "Kendall's p-Value"
$wsInfer.Range("H1").Value = "Spearman's Rho (" + [char]0x03C1 + ")"
$wsInfer.Range("I1").Value = "Spearman's p-Value"

# Reuse the bold/centered/bordered header formatting already used on the
# other sheets instead of re-deriving it property by property.
$wsTypestate.Range("A1:I1").Copy()
$wsInfer.Range("A1:I1").PasteSpecial(-4122)

# Row 2 - "Something - COG Dataset 1"
$wsInfer.Range("A2").Value = "Something - COG Dataset 1"
$wsInfer.Range("B2").Value = 23
$wsInfer.Range("C2").Value = 0
$wsInfer.Range("D2").Value = 0
$wsInfer.Range("E2").Value = "TEMP"
$wsInfer.Range("F2").Value = "TEMP"
$wsInfer.Range("G2").Value = "TEMP"
$wsInfer.Range("H2").Value = "TEMP"
$wsInfer.Range("I2").Value = "TEMP"

# Row 3 - "Human Judgement (readability rating 5 through 1) - COG Dataset 3"
$wsInfer.Range("A3").Value = "Human Judgement (readability rating 5 through 1) - COG Dataset 3"
$wsInfer.Range("B3").Value = 100
$wsInfer.Range("C3").Value = 16
$wsInfer.Range("D3").Value = 16
$wsInfer.Range("E3").Value = 100
$wsInfer.Range("F3").Value = -0.139855453798248
$wsInfer.Range("G3").Value = 0.09054763551591515
$wsInfer.Range("H3").Value = -0.1701042564652381
$wsInfer.Range("I3").Value = 0.09064383040289893

# Row 4 - "Behavioral (correctness in %) - fMRI Dataset"
$wsInfer.Range("A4").Value = "Behavioral (correctness in %) - fMRI Dataset"
$wsInfer.Range("B4").Value = 16
$wsInfer.Range("C4").Value = 0
$wsInfer.Range("D4").Value = 0
$wsInfer.Range("E4").Value = 16

# Row 5 - "Behavioral (time in sec.) - fMRI Dataset"
$wsInfer.Range("A5").Value = "Behavioral (time in sec.) - fMRI Dataset"
$wsInfer.Range("B5").Value = 16
$wsInfer.Range("C5").Value = 0
$wsInfer.Range("D5").Value = 0
$wsInfer.Range("E5").Value = 16

# Row 6 - "Human Judgement (readability low, med, high) - fMRI Dataset"
$wsInfer.Range("A6").Value = "Human Judgement (readability low, med, high) - fMRI Dataset"
$wsInfer.Range("B6").Value = 16
$wsInfer.Range("C6").Value = 0
$wsInfer.Range("D6").Value = 0
$wsInfer.Range("E6").Value = 16

$wsAllTools.Activate()
